# Änderung am 2025-10-13 15:20  auf HelmutsLaptop
# Adds the "Level Korrektur" / EntryNeu explanation block to Tabelle2,
# shrinks the font of the added detail line, tweaks the printed page
# setup for Tabelle2, and updates the active sheet / selections to match
# where the author was last working.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")

# --- New content appended to Tabelle2 (rows 29, 30, 32) ------------------
# Set the cell values first, in the same left-to-right / top-to-bottom
# order the strings were authored, so new shared-string entries land in
# the same order as in the target workbook.
$ws2.Range("A29").Value = "EntryNeu"
$ws2.Range("B29").Value = "Level Korrektur:  "
$ws2.Range("B30").Value = "  alte linien für Peeks;  GAPs, Dojis im Level"
$ws2.Range("B32").Value = "Orders an die Börse"

# Give the new cells the same "wrap, left/top aligned, 14pt" look used by
# the rest of the sheet (copied from an existing cell with that style).
$ws2.Range("B19").Copy()
$ws2.Range("A29").PasteSpecial(-4122)
$ws2.Range("B19").Copy()
$ws2.Range("B29").PasteSpecial(-4122)
$ws2.Range("B19").Copy()
$ws2.Range("A30").PasteSpecial(-4122)
$ws2.Range("B19").Copy()
$ws2.Range("B32").PasteSpecial(-4122)

# The continuation line (B30) uses a smaller, 12pt font.
$ws2.Range("B30").Font.Size = 12
$ws2.Range("B30").WrapText = $true
$ws2.Range("B30").HorizontalAlignment = -4131
$ws2.Range("B30").VerticalAlignment = -4160

# Match the row height the rest of the (14/18pt) sheet uses for these
# newly-populated rows.
$ws2.Rows.Item(29).RowHeight = 18
$ws2.Rows.Item(30).RowHeight = 18
$ws2.Rows.Item(32).RowHeight = 18

# B24 ("alle off Entries unterhalb aus") shrinks from 14pt to 12pt too.
$ws2.Range("B24").Font.Size = 12

# --- Page setup for Tabelle2 ----------------------------------------------
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- Window / selection state ---------------------------------------------
# The author had scrolled down and was working on Tabelle2 when the file
# was saved, with Tabelle1's selection left on A19:B20 and Tabelle2's on
# B24.
$ws1.Range("A19:B20").Select()
$ws2.Select()
$ws2.Range("B24").Select()
